$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (sharedStrings footer text)
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 15:16"

# Country/provincias data refresh: several countries were re-ordered/
# inserted in the source list and case counts updated. Rewrite each
# affected row (country name + the 7 numeric columns) in place so the
# sheet ends up matching the refreshed snapshot.

# Row 30: Pakistan
$ws.Range("A30").Value = "Pakistan"
$ws.Range("B30").Value = 873
$ws.Range("C30").Value = 97
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = 854
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 6

# Row 61: Serbia
$ws.Range("A61").Value = "Serbia"
$ws.Range("B61").Value = 249
$ws.Range("C61").Value = 27
$ws.Range("D61").Value = 2
$ws.Range("E61").Value = 245
$ws.Range("F61").Value = 4
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 2

# Row 62: Republica Dominicana
$ws.Range("A62").Value = "Republica Dominicana"
$ws.Range("B62").Value = 245
$ws.Range("C62").Value = 43
$ws.Range("D62").Value = 0
$ws.Range("E62").Value = 242
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H62").Value = 3

# Row 63: Colombia
$ws.Range("A63").Value = "Colombia"
$ws.Range("B63").Value = 235
$ws.Range("C63").Value = 4
$ws.Range("D63").Value = 3
$ws.Range("E63").Value = 229
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 3

# Row 81: Vietnam
$ws.Range("A81").Value = "Vietnam"
$ws.Range("B81").Value = 123
$ws.Range("C81").Value = 10
$ws.Range("D81").Value = 17
$ws.Range("E81").Value = 106
$ws.Range("F81").Value = 2
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 0

# Row 97: Azerbaiyan
$ws.Range("A97").Value = "Azerbaiyan"
$ws.Range("B97").Value = 72
$ws.Range("C97").Value = 7
$ws.Range("D97").Value = 10
$ws.Range("E97").Value = 61
$ws.Range("F97").Value = 3
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 1

# Row 98: Reunion
$ws.Range("A98").Value = "Reunion"
$ws.Range("B98").Value = 71
$ws.Range("C98").Value = 7
$ws.Range("D98").Value = 1
$ws.Range("E98").Value = 70
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0

# Row 99: Senegal
$ws.Range("A99").Value = "Senegal"
$ws.Range("B99").Value = 67
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 5
$ws.Range("E99").Value = 62
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0

# Row 100: Oman
$ws.Range("A100").Value = "Oman"
$ws.Range("B100").Value = 66
$ws.Range("C100").Value = 11
$ws.Range("D100").Value = 17
$ws.Range("E100").Value = 49
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 0

# Row 111: Consejo Danes para los Refugiados
$ws.Range("A111").Value = "Consejo Danes para los Refugiados"
$ws.Range("B111").Value = 36
$ws.Range("C111").Value = 6
$ws.Range("D111").Value = 0
$ws.Range("E111").Value = 35
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = 1

# Row 112: Nigeria
$ws.Range("A112").Value = "Nigeria"
$ws.Range("B112").Value = 36
$ws.Range("C112").Value = 6
$ws.Range("D112").Value = 2
$ws.Range("E112").Value = 33
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 1

# Row 113: Cuba
$ws.Range("A113").Value = "Cuba"
$ws.Range("B113").Value = 35
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 0
$ws.Range("E113").Value = 34
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 1

# Row 114: Banglades
$ws.Range("A114").Value = "Banglades"
$ws.Range("B114").Value = 33
$ws.Range("C114").Value = 6
$ws.Range("D114").Value = 5
$ws.Range("E114").Value = 25
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 3

# Row 115: Puerto Rico
$ws.Range("A115").Value = "Puerto Rico"
$ws.Range("B115").Value = 31
$ws.Range("C115").Value = 8
$ws.Range("D115").Value = 1
$ws.Range("E115").Value = 28
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 2

# Row 127: Guayana Francesa
$ws.Range("A127").Value = "Guayana Francesa"
$ws.Range("B127").Value = 20
$ws.Range("C127").Value = 2
$ws.Range("D127").Value = 6
$ws.Range("E127").Value = 14
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 0

# Row 128: Ruanda
$ws.Range("A128").Value = "Ruanda"
$ws.Range("B128").Value = 19
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 0
$ws.Range("E128").Value = 19
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

# Row 129: Guyana
$ws.Range("A129").Value = "Guyana"
$ws.Range("B129").Value = 19
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 18
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

# Row 130: Guatemala
$ws.Range("A130").Value = "Guatemala"
$ws.Range("B130").Value = 19
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 0
$ws.Range("E130").Value = 18
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 1

# Row 131: Jamaica
$ws.Range("A131").Value = "Jamaica"
$ws.Range("B131").Value = 19
$ws.Range("C131").Value = 0
$ws.Range("D131").Value = 2
$ws.Range("E131").Value = 16
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 1

# Row 132: Togo
$ws.Range("A132").Value = "Togo"
$ws.Range("B132").Value = 18
$ws.Range("C132").Value = 2
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 18
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 0

# Row 133: Polinesia Francesa
$ws.Range("A133").Value = "Polinesia Francesa"
$ws.Range("B133").Value = 18
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 0
$ws.Range("E133").Value = 18
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

# Row 137: Gibraltar
$ws.Range("A137").Value = "Gibraltar"
$ws.Range("B137").Value = 15
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 5
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

# Row 143: Guinea Ecuatorial
$ws.Range("A143").Value = "Guinea Ecuatorial"
$ws.Range("B143").Value = 9
$ws.Range("C143").Value = 3
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 9
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

# Row 144: Aruba
$ws.Range("A144").Value = "Aruba"
$ws.Range("B144").Value = 9
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 1
$ws.Range("E144").Value = 8
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

# Row 145: Nueva Caledonia
$ws.Range("A145").Value = "Nueva Caledonia"
$ws.Range("B145").Value = 8
$ws.Range("C145").Value = 4
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 8
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

# Row 146: Seychelles
$ws.Range("A146").Value = "Seychelles"
$ws.Range("B146").Value = 7
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 0
$ws.Range("E146").Value = 7
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 0

# Row 150: Surinam
$ws.Range("A150").Value = "Surinam"
$ws.Range("B150").Value = 5
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 0
$ws.Range("E150").Value = 5
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 0

# Row 151: Haiti
$ws.Range("A151").Value = "Haiti"
$ws.Range("B151").Value = 5
$ws.Range("C151").Value = 3
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 5
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0
